# Replace the final "How to Make a Wireless Bicycle Outfit..." image-caption
# paragraph (a leftover Markdown-style figure placeholder, spread across
# several runs) with a single typeset editorial placeholder run:
#   [INSERT FIGURE 9.1 NEAR HERE]

$d = $word.ActiveDocument

# Locate the target paragraph robustly by searching for a short, unique
# substring of its text rather than assuming it is literally the last
# paragraph in the document.
$anchor = $d.Content.Duplicate
$found = $anchor.Find.Execute("submitted to the September 1910 issue", $true, $false, `
                               $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the bicycle-wireless figure caption paragraph."
}

$targetPara = $anchor.Paragraphs(1)
$targetRange = $targetPara.Range

# Use InsertXML (rather than a plain Find/Replace or Range.Text assignment)
# so the resulting run keeps xml:space="preserve" on <w:t>, matching this
# document's run-serialization convention, and so the whole paragraph's
# multi-run content is atomically replaced by exactly one run.
$newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
                    '<w:r><w:t xml:space="preserve">[INSERT FIGURE 9.1 NEAR HERE]</w:t></w:r>' + `
                    '</w:p>'

$targetRange.InsertXML($newParagraphXml)
